$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "261.32"
$ws.Range("D2").Style = "Normal"

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "22.75"
$ws.Range("D3").Style = "Normal"

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "6.199"
$ws.Range("D4").Style = "Normal"

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.06121"
$ws.Range("D5").Style = "Normal"

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "6.741"
$ws.Range("D6").Style = "Normal"

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.499"
$ws.Range("D7").Style = "Normal"

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.7981"
$ws.Range("D9").Style = "Normal"

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1578"
$ws.Range("D10").Style = "Normal"

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.08086"
$ws.Range("D11").Style = "Normal"

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.03315"
$ws.Range("D12").Style = "Normal"

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.03107"
$ws.Range("D13").Style = "Normal"

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.09290"
$ws.Range("D14").Style = "Normal"

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.915"
$ws.Range("D15").Style = "Normal"

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.001700"
$ws.Range("D16").Style = "Normal"

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.04842"
$ws.Range("D17").Style = "Normal"

$ws.Range("B18").Value = "TigerCash"

$ws.Range("C18").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.006201"
$ws.Range("D18").Style = "Normal"

$ws.Range("E18").Value = "17TigerCashTCH"

$ws.Range("B19").Value = "BitKan"

$ws.Range("C19").Value = "https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan"

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.001099"
$ws.Range("D19").Style = "Normal"

$ws.Range("E19").Value = "18BitKanKAN"

$ws.Range("B20").Value = "HotbitToken"

$ws.Range("C20").Value = "https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb"

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.003398"
$ws.Range("D20").Style = "Normal"

$ws.Range("E20").Value = "19HotbitTokenHTB"

$ws.Range("B21").Value = "NitroEx"

$ws.Range("C21").Value = "https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx"

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.0001499"
$ws.Range("D21").Style = "Normal"

$ws.Range("E21").Value = "20NitroExNTX"

$ws.Range("B22").Value = "LEO"

$ws.Range("C22").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "3.693"
$ws.Range("D22").Style = "Normal"

$ws.Range("E22").Value = "21LEOLEO"

$ws.Range("B23").Value = "BTSEToken"

$ws.Range("C23").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "2.259"
$ws.Range("D23").Style = "Normal"

$ws.Range("E23").Value = "22BTSETokenBTSE"

$ws.Range("B24").Value = "One"

$ws.Range("C24").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.01332"
$ws.Range("D24").Style = "Normal"

$ws.Range("E24").Value = "23OneONE"

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.3361"
$ws.Range("D25").Style = "Normal"

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0006163"
$ws.Range("D27").Style = "Normal"

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.04587"
$ws.Range("D40").Style = "Normal"

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007157"
$ws.Range("D41").Style = "Normal"

$ws.Range("B42").Value = "CEJI"

$ws.Range("C42").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.003899"
$ws.Range("D42").Style = "Normal"

$ws.Range("E42").Value = "41CEJICEJI"

$ws.Range("B43").Value = "BKEXToken"

$ws.Range("C43").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1121"
$ws.Range("D43").Style = "Normal"

$ws.Range("E43").Value = "42BKEXTokenBKK"

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.01022"
$ws.Range("D44").Style = "Normal"

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00006018"
$ws.Range("D46").Style = "Normal"

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.7498"
$ws.Range("D48").Style = "Normal"

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.06184"
$ws.Range("D49").Style = "Normal"

$ws.Range("E49").Value = "48BOLOBOLOWorstin24h"

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.00002099"
$ws.Range("D50").Style = "Normal"
